# Acetic Acid Area_output.xlsx — "hoac pt(100) figs and code updates"
#
# The single data sheet is rebuilt: an extra "L" column moves from J to A
# (gaining the same bold/bordered/centered header style, s=1, already used
# by the other header cells and by the original A2), the species-name
# headers are replaced by a new list (Water, Acetic acid43, Acetaldehyde,
# CO2, CO, Ketene, H2, Acetic acid60), the old two-row TPD-area table is
# replaced by a new 9-row table (rows 2-10) of -1 placeholders plus CO/H2
# peak-area figures, and the now-unused column J is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- formatting first, while the "donor" styled cells (B1, A2) still exist ---

# A1 is a brand-new cell; give it the same header style (bold, thin border,
# centered) already carried by B1:I1 so styles.xml gains no new cellXfs.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# A3:A10 are brand-new cells in column A; give them the same style already
# carried by A2 (bold, thin border, centered).
$ws.Range("A2").Copy()
$ws.Range("A3:A10").PasteSpecial(-4122)

# Column J no longer exists in the new layout.
$ws.Range("J1:J2").Clear()

# --- header row (A1:I1) ---

$headers = @("L","Water","Acetic acid43","Acetaldehyde","CO2","CO","Ketene","H2","Acetic acid60")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- data rows (A2:I10) ---
# Columns: A=L, B=Water, C=Acetic acid43, D=Acetaldehyde, E=CO2, F=CO,
#          G=Ketene, H=H2, I=Acetic acid60

$data = @(
    @(0,      -1, -1, -1, -1, [double]"6.023226644740085e-07", -1, [double]"1.344941710783223e-06", -1),
    @(0,      -1, -1, -1, -1, [double]"8.662993879037748e-07", -1, [double]"1.424802033287694e-06", -1),
    @(0.0004, -1, -1, -1, -1, [double]"2.79860466923858e-07",  -1, [double]"1.071084382567251e-06", -1),
    @(0.0008, -1, -1, -1, -1, [double]"5.329421183453525e-07", -1, [double]"9.715907029045017e-07", -1),
    @(0.0012, -1, -1, -1, -1, [double]"5.73752101054891e-07",  -1, [double]"8.498520080187098e-07", -1),
    @(0.0015, -1, -1, -1, -1, [double]"7.416994930621231e-07", -1, [double]"1.291435777499689e-06", -1),
    @(0.003,  -1, -1, -1, -1, [double]"1.20699759654942e-06",  -1, [double]"2.102431208430838e-06", -1),
    @(0.004,  -1, -1, -1, -1, [double]"7.401729158066803e-07", -1, [double]"9.731101506668055e-07", -1),
    @(0.005,  -1, -1, -1, -1, [double]"2.924680721579058e-06", -1, [double]"2.263300604122058e-06", -1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
